$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 31 (shifts rows 31..123 down to 32..124,
# copying formatting from the row above, matching Excel's native Insert behaviour).
$ws.Rows("31:31").Insert()

# Populate the newly inserted row 31 with this week's new data point.
$ws.Range("A31").Value = 4
$ws.Range("B31").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C31").Value = "Los Lagos"
$ws.Range("D31").Value = 44487
$ws.Range("E31").Value = 10
$ws.Range("F31").Value = 100112039
$ws.Range("G31").Value = "Ciboulette"
$ws.Range("H31").Value = "Sin especificar"
$ws.Range("I31").Value = "Primera"
$ws.Range("J31").Value = 80
$ws.Range("K31").Value = 2500
$ws.Range("L31").Value = 2500
$ws.Range("M31").Value = 2500
$ws.Range("N31").Value = '$/docena de atados'
$ws.Range("O31").Value = "Región Metropolitana"
$ws.Range("P31").Value = 833
$ws.Range("Q31").Value = 3
$ws.Range("R31").Value = "Hortaliza"
